$wb = $excel.ActiveWorkbook

# Sheet1: ip_address_list - add row 2
$ws1 = $wb.Worksheets.Item("ip_address_list")
$ws1.Range("A2").Value = "518"
$ws1.Range("B2").Value = "192.168.1.241"
$ws1.Range("C2").Value = "255.255.255.0"
$ws1.Range("D2").Value = "pozngg"
$ws1.Range("E2").Value = 0

# Sheet3: disk_list - update A1, C1 and add row 4
$ws3 = $wb.Worksheets.Item("disk_list")
$ws3.Range("A1").Value = "518-2"
$ws3.Range("C1").Value = "\\192.168.208.200"
$ws3.Range("A4").Value = "518"
$ws3.Range("B4").Value = "V"
$ws3.Range("C4").Value = "\\192.168.1.10\10_vision"
$ws3.Range("D4").Value = "jhv_vision"
$ws3.Range("E4").Value = "Jhv*2708"
$ws3.Range("F4").Value = "druhá síť Valeo"

# Sheet4: Settings - update values
$ws4 = $wb.Worksheets.Item("Settings")
$ws4.Range("B3").Value = 0
$ws4.Range("B4").Value = 0
$ws4.Range("A6").Value = "aktualizovat statusy disku (default)"
$ws4.Range("B6").Value = 0
$ws4.Range("B6").Select()
